$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the mixed-up fixed interactions in rows 2 and 3 -----------------
# Row 2 should describe EGF->EGFR, row 3 should describe TGFb->TGFRb
# (their "From"/"To" values had been swapped). Read the old values first.
$a2 = $ws.Range("A2").Value2
$c2 = $ws.Range("C2").Value2
$a3 = $ws.Range("A3").Value2
$c3 = $ws.Range("C3").Value2

$ws.Range("A2").Value2 = $a3
$ws.Range("C2").Value2 = $c3
$ws.Range("A3").Value2 = $a2
$ws.Range("C3").Value2 = $c2

# --- Re-key the D-column concatenation formulas ---------------------------
# D2 and D3 keep their own (non-shared) formulas; re-enter them so the
# computed display values refresh for the swapped A/C values above.
$ws.Range("D2").Formula = "=A2&B2&C2"
$ws.Range("D3").Formula = "=A3&B3&C3"
# D4:D44 share one formula (anchored at D4 instead of D3 as before).
$ws.Range("D4:D44").Formula = "=A4&B4&C4"

# --- Restore the saved view/selection state --------------------------------
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("A2:F2").Select()
